$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Fix up sheet1 ("VIP Promos") selection / active-view state
# ------------------------------------------------------------------
$vip = $wb.Worksheets.Item("VIP Promos")
$vip.Range("E14").Select()

# ------------------------------------------------------------------
# 2. Add the new "20 promos" sheet after "Staging"
# ------------------------------------------------------------------
$staging = $wb.Worksheets.Item("Staging")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $staging)
$newSheet.Name = "20 promos"

# Header row (row 1) - same headers as "Staging"
$headers = @("Name","Created By Email","Created By Name","Discount","Created at","Starts","Ends","Active week days","Active time","Available types","Available locations","Amount of uses","Amount of Subscription Valid Periods","Max Active at Same Time (Quantity)","Max Active at Same Time (Grace Period)","Min Hours","Max Hours","Once per account","Only Valid Before Parking Starts","Notes","Business Account ID","Business Account Name","Batch Name","Fee Amount","Event Exempt","Refund Service Fee","No Discount Or Validation Restrictions","Tag")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $newSheet.Cells.Item(1, $col).Value = $headers[$i]
}
$newSheet.Cells.Item(1, 15).NumberFormat = "0"
$newSheet.Cells.Item(1, 21).NumberFormat = "0"
$newSheet.Cells.Item(1, 24).NumberFormat = '\$#,##0.00_);\(\$#,##0.00\)'

# Promo code (col A) for each of the 20 rows
$promoCodes = @("MS101VIP","DA102VIP","JK103VIP","BL104VIP","JP105VIP","CW106VIP","DM107VIP","BS108VIP","BL109VIP","KP110VIP","HC111VIP","CM112VIP","GT113VIP","WM114VIP","SG115VIP","DL116VIP","JC117VIP","MC118VIP","LD119VIP","TG120VIP")

for ($r = 0; $r -lt $promoCodes.Length; $r++) {
    $row = $r + 2

    $newSheet.Cells.Item($row, 1).Value = $promoCodes[$r]                         # A Name
    $newSheet.Cells.Item($row, 2).Value = "sdumas@premiumparking.com"             # B Created By Email
    $newSheet.Cells.Item($row, 3).Value = "Sarah Dumas"                           # C Created By Name
    $newSheet.Cells.Item($row, 4).NumberFormat = "0%"                             # D Discount
    $newSheet.Cells.Item($row, 4).Value = 1
    $newSheet.Cells.Item($row, 6).Value = "01/01/2023 12:00 AM (CST)"             # F Starts
    $newSheet.Cells.Item($row, 7).Font.Color = 0                                  # G Ends
    $newSheet.Cells.Item($row, 7).Value = "12/31/2023 12:00 AM (CST)"
    $newSheet.Cells.Item($row, 8).Value = "Sunday, Monday, Tuesday, Wednesday, Thursday, Friday, Saturday"   # H
    $newSheet.Cells.Item($row, 9).Value = "All day"                               # I
    $newSheet.Cells.Item($row, 10).Value = "session, reservation"                 # J
    $newSheet.Cells.Item($row, 11).Value = "All locations"                        # K
    $newSheet.Cells.Item($row, 12).Value = 100                                    # L Amount of uses
    $newSheet.Cells.Item($row, 13).Value = "Unlimited"                            # M
    $newSheet.Cells.Item($row, 14).Value = "Unlimited"                            # N
    $newSheet.Cells.Item($row, 15).NumberFormat = "0"                             # O Max Active (Grace Period) - blank
    $newSheet.Cells.Item($row, 16).Value = "Unlimited"                            # P Min Hours
    $newSheet.Cells.Item($row, 17).Value = 24                                     # Q Max Hours
    $newSheet.Cells.Item($row, 18).Value = "No"                                   # R Once per account
    $newSheet.Cells.Item($row, 19).Value = "No"                                   # S Only Valid Before Parking Starts
    $newSheet.Cells.Item($row, 20).Value = "PremiumVIP 2023"                      # T Notes
    $newSheet.Cells.Item($row, 21).NumberFormat = "0"                             # U Business Account ID - blank
    $newSheet.Cells.Item($row, 22).Value = "VIP"                                  # V Business Account Name
    $newSheet.Cells.Item($row, 24).NumberFormat = '\$#,##0.00_);\(\$#,##0.00\)'   # X Fee Amount
    $newSheet.Cells.Item($row, 24).Value = 0
    $newSheet.Cells.Item($row, 25).Value = "No"                                   # Y Event Exempt
    $newSheet.Cells.Item($row, 26).Value = "No"                                   # Z Refund Service Fee
    $newSheet.Cells.Item($row, 27).Value = "No"                                   # AA No Discount Or Validation Restrictions
}

# ------------------------------------------------------------------
# 3. Hyperlinks on column B (Created By Email)
# ------------------------------------------------------------------
$newSheet.Hyperlinks.Add($newSheet.Cells.Item(2, 2), "mailto:sdumas@premiumparking.com")
$rangeRest = $newSheet.Range($newSheet.Cells.Item(3, 2), $newSheet.Cells.Item(21, 2))
$newSheet.Hyperlinks.Add($rangeRest, "mailto:sdumas@premiumparking.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "sdumas@premiumparking.com")

# Re-apply the Hyperlink cell style (Hyperlinks.Add creates a near-duplicate style) so
# every linked cell in column B points at the canonical "Hyperlink" style.
$newSheet.Range($newSheet.Cells.Item(2, 2), $newSheet.Cells.Item(21, 2)).Style = "Hyperlink"

# ------------------------------------------------------------------
# 4. Sheet view state for the new sheet, then make it the active tab
# ------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("H24").Select()
